# "Update comments, sheet names"
#
# - Rename the "Climate" sheet to "Application climate".
# - Remove the two outdated comments on the "Storage EFs" sheet (B1 and D1),
#   keeping the D4 comment as-is.
# - Update the remembered selections: "Slurry & application" moves from
#   D17 to D16, and "Storage EFs" moves from E8 to D10. Selecting the
#   "Slurry & application" range last restores it as the active sheet/tab
#   (matching the workbook's activeTab going back to the first sheet).

$wb = $excel.ActiveWorkbook

# 1. Rename "Climate" -> "Application climate"
$wsClimate = $wb.Worksheets.Item("Climate")
$wsClimate.Name = "Application climate"

# 2. Drop the two resolved/obsolete comments on "Storage EFs", leave D4 alone
$wsStorage = $wb.Worksheets.Item("Storage EFs")
$wsStorage.Range("B1").Comment.Delete()
$wsStorage.Range("D1").Comment.Delete()

# 3. Update the saved selection on "Storage EFs" to D10
$wsStorage.Range("D10").Select()

# 4. Update the saved selection on "Slurry & application" to D16 and make
#    it the active sheet/tab again
$wsSlurry = $wb.Worksheets.Item("Slurry & application")
$wsSlurry.Range("D16").Select()
